$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.025010761346032
$ws.Range("D2").Value = 1.029673052772681
$ws.Range("E2").Value = 1.028637437892254
$ws.Range("F2").Value = 1.034984257447097
$ws.Range("I2").Value = 1.030812905386939
$ws.Range("J2").Value = 1.03018217148517
$ws.Range("K2").Value = 1.032486282148583
$ws.Range("L2").Value = 1.03145367510015
$ws.Range("M2").Value = 1.037782163779857
$ws.Range("N2").Value = 1.03164514833223

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.026123957834818
$ws.Range("D3").Value = 1.030498944535342
$ws.Range("E3").Value = 1.029698296874913
$ws.Range("F3").Value = 1.036280304298826
$ws.Range("I3").Value = 1.031024796345531
$ws.Range("J3").Value = 1.03093367399596
$ws.Range("K3").Value = 1.033120320915411
$ws.Range("L3").Value = 1.032321829617442
$ws.Range("M3").Value = 1.03888621682905
$ws.Range("N3").Value = 1.03239771806278

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.026844189904897
$ws.Range("D4").Value = 1.031033052556245
$ws.Range("E4").Value = 1.030385034363529
$ws.Range("F4").Value = 1.037118987113625
$ws.Range("I4").Value = 1.031160423402204
$ws.Range("J4").Value = 1.031419348842769
$ws.Range("K4").Value = 1.033529643392231
$ws.Range("L4").Value = 1.032883282161468
$ws.Range("M4").Value = 1.039600125297586
$ws.Range("N4").Value = 1.032884082623579

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.02714695734478
$ws.Range("D5").Value = 1.031257519973514
$ws.Range("E5").Value = 1.030673808987488
$ws.Range("F5").Value = 1.0374715841988
$ws.Range("I5").Value = 1.031217086615875
$ws.Range("J5").Value = 1.031623383981492
$ws.Range("K5").Value = 1.033701497015113
$ws.Range("L5").Value = 1.033119245153202
$ws.Range("M5").Value = 1.039900137247266
$ws.Range("N5").Value = 1.033088407515602

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.027197792300881
$ws.Range("D6").Value = 1.031295204842994
$ws.Range("E6").Value = 1.030722299649635
$ws.Range("F6").Value = 1.03753078772697
$ws.Range("I6").Value = 1.031226579837942
$ws.Range("J6").Value = 1.031657634032904
$ws.Range("K6").Value = 1.03373033877451
$ws.Range("L6").Value = 1.033158860190772
$ws.Range("M6").Value = 1.039950503866045
$ws.Range("N6").Value = 1.033122706206015

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.026848235568017
$ws.Range("D7").Value = 1.031036052182155
$ws.Range("E7").Value = 1.030388892707706
$ws.Range("F7").Value = 1.03712369847244
$ws.Range("I7").Value = 1.031161181931334
$ws.Range("J7").Value = 1.031422075729872
$ws.Range("K7").Value = 1.033531940594741
$ws.Range("L7").Value = 1.032886435391437
$ws.Range("M7").Value = 1.039604134522514
$ws.Range("N7").Value = 1.032886813383174

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.025386988669831
$ws.Range("D8").Value = 1.029952228725643
$ws.Range("E8").Value = 1.028995900285317
$ws.Range("F8").Value = 1.03542225256963
$ws.Range("I8").Value = 1.03088482151839
$ws.Range("J8").Value = 1.030436269290988
$ws.Range("K8").Value = 1.032700753900911
$ws.Range("L8").Value = 1.031747134719611
$ws.Range("M8").Value = 1.03815538560631
$ws.Range("N8").Value = 1.031899606986074

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.022811416673853
$ws.Range("D9").Value = 1.028040099978976
$ws.Range("E9").Value = 1.02654347076985
$ws.Range("F9").Value = 1.032424412683032
$ws.Range("I9").Value = 1.030386497349243
$ws.Range("J9").Value = 1.028694553941567
$ws.Range("K9").Value = 1.031228864506197
$ws.Range("L9").Value = 1.029737208726809
$ws.Range("M9").Value = 1.035598698774354
$ws.Range("N9").Value = 1.030155418201078

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.021093844462891
$ws.Range("D10").Value = 1.026763799670111
$ws.Range("E10").Value = 1.024909958425509
$ws.Range("F10").Value = 1.030425947166272
$ws.Range("I10").Value = 1.030046650855593
$ws.Range("J10").Value = 1.0275302853801
$ws.Range("K10").Value = 1.03024272595997
$ws.Range("L10").Value = 1.028395654055341
$ws.Range("M10").Value = 1.03389158010796
$ws.Range("N10").Value = 1.028989496244711

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.020349974927385
$ws.Range("D11").Value = 1.026210778149663
$ws.Range("E11").Value = 1.024202961615616
$ws.Range("F11").Value = 1.029560583534112
$ws.Range("I11").Value = 1.029897682249558
$ws.Range("J11").Value = 1.027025394884486
$ws.Range("K11").Value = 1.029814555434935
$ws.Range("L11").Value = 1.027814356872064
$ws.Range("M11").Value = 1.033151724393654
$ws.Range("N11").Value = 1.028483888746681

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.02007364460503
$ws.Range("D12").Value = 1.026005304679156
$ws.Range("E12").Value = 1.023940399388156
$ws.Range("F12").Value = 1.029239144172345
$ws.Range("I12").Value = 1.029842076125353
$ws.Range("J12").Value = 1.026837741920234
$ws.Range("K12").Value = 1.029655338176648
$ws.Range("L12").Value = 1.027598376996463
$ws.Range("M12").Value = 1.032876807456868
$ws.Range("N12").Value = 1.028295969293696

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.020132919511595
$ws.Range("D13").Value = 1.026049382004774
$ws.Range("E13").Value = 1.023996717726187
$ws.Range("F13").Value = 1.029308094263726
$ws.Range("I13").Value = 1.029854016164388
$ws.Range("J13").Value = 1.026877999296578
$ws.Range("K13").Value = 1.029689498775783
$ws.Range("L13").Value = 1.027644708151711
$ws.Range("M13").Value = 1.032935782703851
$ws.Range("N13").Value = 1.028336283840132

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.020327133876767
$ws.Range("D14").Value = 1.026193794806015
$ws.Range("E14").Value = 1.024181257168929
$ws.Range("F14").Value = 1.029534013368402
$ws.Range("I14").Value = 1.029893091391069
$ws.Range("J14").Value = 1.027009885765238
$ws.Range("K14").Value = 1.029801398072911
$ws.Range("L14").Value = 1.027796505141568
$ws.Range("M14").Value = 1.033129001758393
$ws.Range("N14").Value = 1.028468357602705

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.020446792521311
$ws.Range("D15").Value = 1.026282764782581
$ws.Range("E15").Value = 1.024294964325012
$ws.Range("F15").Value = 1.029673208875705
$ws.Range("I15").Value = 1.029917130806397
$ws.Range("J15").Value = 1.027091130221042
$ws.Range("K15").Value = 1.029870319631845
$ws.Range("L15").Value = 1.027890024271757
$ws.Range("M15").Value = 1.033248036870955
$ws.Range("N15").Value = 1.028549717434955

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.021143209348719
$ws.Range("D16").Value = 1.026800493958978
$ws.Range("E16").Value = 1.024956886239272
$ws.Range("F16").Value = 1.030483377970121
$ws.Range("I16").Value = 1.030056499200577
$ws.Range("J16").Value = 1.027563777334712
$ws.Range("K16").Value = 1.030271117619957
$ws.Range("L16").Value = 1.028434224436823
$ws.Range("M16").Value = 1.0339406677657
$ws.Range("N16").Value = 1.02902303576174

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.021580012085836
$ws.Range("D17").Value = 1.027125151164205
$ws.Range("E17").Value = 1.025372178673505
$ws.Range("F17").Value = 1.030991569862687
$ws.Range("I17").Value = 1.030143435829492
$ws.Range("J17").Value = 1.027860053814259
$ws.Range("K17").Value = 1.030522215094061
$ws.Range("L17").Value = 1.028775480428284
$ws.Range("M17").Value = 1.034374958082238
$ws.Range("N17").Value = 1.029319732987872

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.021834777421919
$ws.Range("D18").Value = 1.027314482051774
$ws.Range("E18").Value = 1.025614443126939
$ws.Range("F18").Value = 1.031287988331303
$ws.Range("I18").Value = 1.030193969595375
$ws.Range("J18").Value = 1.028032794124933
$ws.Range("K18").Value = 1.030668563530098
$ws.Range("L18").Value = 1.028974491204859
$ws.Range("M18").Value = 1.034628208538501
$ws.Range("N18").Value = 1.0294927186096

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.021921643439301
$ws.Range("D19").Value = 1.027379032866259
$ws.Range("E19").Value = 1.02569705439965
$ws.Range("F19").Value = 1.0313890592318
$ws.Range("I19").Value = 1.030211170643616
$ws.Range("J19").Value = 1.0280916817918
$ws.Range("K19").Value = 1.030718445502041
$ws.Range("L19").Value = 1.029042342305345
$ws.Range("M19").Value = 1.034714549652737
$ws.Range("N19").Value = 1.029551689903709

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.021533148777545
$ws.Range("D20").Value = 1.027090322261634
$ws.Range("E20").Value = 1.02532761848401
$ws.Range("F20").Value = 1.030937045830918
$ws.Range("I20").Value = 1.030134126443914
$ws.Range("J20").Value = 1.027828273682101
$ws.Range("K20").Value = 1.030495286348311
$ws.Range("L20").Value = 1.028738870833708
$ws.Range("M20").Value = 1.034328369449033
$ws.Range("N20").Value = 1.02928790772428

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.020269943263861
$ws.Range("D21").Value = 1.026151270401866
$ws.Range("E21").Value = 1.024126913588722
$ws.Range("F21").Value = 1.029467485975322
$ws.Range("I21").Value = 1.029881592236041
$ws.Range("J21").Value = 1.026971051645011
$ws.Range("K21").Value = 1.029768451371889
$ws.Range("L21").Value = 1.027751806380705
$ws.Range("M21").Value = 1.033072106383863
$ws.Range("N21").Value = 1.028429468333572

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.019475575843263
$ws.Range("D22").Value = 1.025560523285968
$ws.Range("E22").Value = 1.023372259242807
$ws.Range("F22").Value = 1.028543486457861
$ws.Range("I22").Value = 1.029721236963201
$ws.Range("J22").Value = 1.026431421271724
$ws.Range("K22").Value = 1.029310444814105
$ws.Range("L22").Value = 1.027130850947126
$ws.Range("M22").Value = 1.032281655159392
$ws.Range("N22").Value = 1.027889071623259

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.019896698848025
$ws.Range("D23").Value = 1.025873720697921
$ws.Range("E23").Value = 1.023772289895949
$ws.Range("F23").Value = 1.029033319675924
$ws.Range("I23").Value = 1.029806393873114
$ws.Range("J23").Value = 1.026717552450714
$ws.Range("K23").Value = 1.029553339275255
$ws.Range("L23").Value = 1.02746006452233
$ws.Range("M23").Value = 1.032700744825745
$ws.Range("N23").Value = 1.028175609141343

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.021554324333112
$ws.Range("D24").Value = 1.027106060055066
$ws.Range("E24").Value = 1.025347753218218
$ws.Range("F24").Value = 1.03096168289715
$ws.Range("I24").Value = 1.030138333494843
$ws.Range("J24").Value = 1.027842633980359
$ws.Range("K24").Value = 1.030507454636685
$ws.Range("L24").Value = 1.028755413250985
$ws.Range("M24").Value = 1.034349421043351
$ws.Range("N24").Value = 1.02930228841581

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.023477350669074
$ws.Range("D25").Value = 1.028534703451409
$ws.Range("E25").Value = 1.027177225766796
$ws.Range("F25").Value = 1.033199400394412
$ws.Range("I25").Value = 1.03051667064221
$ws.Range("J25").Value = 1.029145377459356
$ws.Range("K25").Value = 1.031610241926529
$ws.Range("L25").Value = 1.030257102863308
$ws.Range("M25").Value = 1.036260125115098
$ws.Range("N25").Value = 1.030606881939973
